$d = $word.ActiveDocument

# Locate the paragraph that starts the "Part 7" git-walkthrough block
# ("I submitted the following commands ...").  The whole block of
# command/output paragraphs that follows it, all the way through the
# last (bookmark-holding) empty paragraph at the very end of the
# document, gets collapsed into one short summary paragraph plus a
# single trailing empty paragraph.  The _GoBack bookmark is carried
# over onto the end of the new summary paragraph.

$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "I submitted the following commands*") {
        $startPara = $p
        break
    }
}

if ($startPara -eq $null) {
    throw "Could not locate the 'I submitted the following commands' paragraph"
}

$endPara = $d.Paragraphs.Last

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>I forked the courses repo, cloned it on my PC, updated the README with my information, added/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>commited</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the changes and then created a pull request.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>
'@

$r.InsertXML($xml) | Out-Null
